# Final version of the project:
# Move the "Математика" / "Казанский Университет Вычислений;" data row
# (previously the last data row, row 5) up to become the first data row
# (row 2), shifting the other rows (Медицина, Физика, Лингвистика) down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2: Математика
$ws.Range("A2").Value = "Математика"
$ws.Range("B2").Value = 0.0
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 1.0
$ws.Range("E2").Value = "Казанский Университет Вычислений;"

# New row 3: Медицина
$ws.Range("A3").Value = "Медицина"
$ws.Range("B3").Value = 4.329999923706055
$ws.Range("C3").Value = 3.0
$ws.Range("D3").Value = 3.0
$ws.Range("E3").Value = "Самарский Медицинский Институт;Московский Государственный Медицинский Университет;Тамбовский Университет Медицины;"

# New row 4: Физика
$ws.Range("A4").Value = "Физика"
$ws.Range("B4").Value = 4.539999961853027
$ws.Range("C4").Value = 8.0
$ws.Range("D4").Value = 2.0
$ws.Range("E4").Value = "Московский Придуманный Институт;Московский Выдуманный Университет;"

# New row 5: Лингвистика
$ws.Range("A5").Value = "Лингвистика"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("D5").Value = 1.0
$ws.Range("E5").Value = "Воронежский Литературно-Переводческий Университет;"

$wb.Save()
